$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties) right after the
# existing "Unnamed: 28" header in AC1, reusing that cell's formatting
# (bold, centered, bordered) via copy/paste-format so no new style is
# minted.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team's season record (89 wins, 73 losses, 0 ties) for every
# player row.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 73
    $ws.Cells.Item($r, 32).Value = 0
}
